# Feat: Shield and Potion test
#
# This script reproduces, via the Excel COM object model, the edit that:
#  - inserts two new columns ("parent" and "textures") after column A
#  - adds two new rows (Copper Shield / item/shield, Custom Potion / item/potion)
#  - moves the conditional formatting + hyperlinks to their new column positions
#  - updates the selected cell

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Insert two blank columns before the old column B (shifts B..G -> D..I)
# ------------------------------------------------------------------
$ws.Range("B:C").Insert()

# ------------------------------------------------------------------
# 2. Header row (row 1)
# ------------------------------------------------------------------
$ws.Range("B1").Value = "parent"
$ws.Range("C1").Value = "textures"

# ------------------------------------------------------------------
# 3. Fill in the new "parent" / "textures" columns for the existing rows
# ------------------------------------------------------------------
$ws.Range("B2").Value = "item/generated"
$ws.Range("B3").Value = "item/generated"

$ws.Range("C2").Formula = '=_xlfn.CONCAT("{''layer0'':''",A2,"''}")'
$ws.Range("C3").Formula = '=_xlfn.CONCAT("{''layer0'':''",A3,"''}")'

# give the new B2:C3 cells the same (border-less) look as column A in those rows
$ws.Range("A2").Copy($ws.Range("B2:C3"))
$ws.Range("C2").Formula = '=_xlfn.CONCAT("{''layer0'':''",A2,"''}")'
$ws.Range("C3").Formula = '=_xlfn.CONCAT("{''layer0'':''",A3,"''}")'
$ws.Range("B2").Value = "item/generated"
$ws.Range("B3").Value = "item/generated"

# ------------------------------------------------------------------
# 4. Update the command formulas so they reference the shifted
#    custom_model_data column (old $B -> new $D)
# ------------------------------------------------------------------
$ws.Range("I2").Formula = '=_xlfn.CONCAT("/give @s minecraft:", RIGHT($A2,LEN($A2)-FIND("/",$A2)), "{CustomModelData:",$D2,"}")'
$ws.Range("I3").Formula = '=_xlfn.CONCAT("/give @s minecraft:", RIGHT($A3,LEN($A3)-FIND("/",$A3)), "{CustomModelData:",$D3,"}")'

# ------------------------------------------------------------------
# 5. New row 4 : Copper Shield
# ------------------------------------------------------------------
$ws.Range("A4").Value = "item/shield"
$ws.Range("B4").Value = "item/generated"
$ws.Range("C4").Formula = '=_xlfn.CONCAT("{''layer0'':''",A4,"''}")'
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = "blockgame:shield/shield_copper"
$ws.Range("F4").Value = "Copper Shield"
$ws.Range("I4").Formula = '=_xlfn.CONCAT("/give @s minecraft:", RIGHT($A4,LEN($A4)-FIND("/",$A4)), "{CustomModelData:",$D4,"}")'

# match styling: A4/D4/E4/I4 plain, B4/C4 like column A data cells, F4 like column F data cells
$ws.Range("A2").Copy($ws.Range("B4:C4"))
$ws.Range("C4").Formula = '=_xlfn.CONCAT("{''layer0'':''",A4,"''}")'
$ws.Range("B4").Value = "item/generated"
$ws.Range("F2").Copy($ws.Range("F4"))
$ws.Range("F4").Value = "Copper Shield"

# ------------------------------------------------------------------
# 6. New row 5 : Custom Potion
# ------------------------------------------------------------------
$ws.Range("A5").Value = "item/potion"
$ws.Range("B5").Value = "item/generated"
$ws.Range("C5").Value = "{'layer0':'item/potion_overlay', 'layer1':'item/potion'}"
$ws.Range("D5").Value = 4
$ws.Range("E5").Value = "blockgame:potion/potion_custom"
$ws.Range("F5").Value = "Custom Potion"
$ws.Range("I5").Formula = '=_xlfn.CONCAT("/give @s minecraft:", RIGHT($A5,LEN($A5)-FIND("/",$A5)), "{CustomModelData:",$D5,"}")'

$ws.Range("A2").Copy($ws.Range("B5:C5"))
$ws.Range("C5").Value = "{'layer0':'item/potion_overlay', 'layer1':'item/potion'}"
$ws.Range("B5").Value = "item/generated"
$ws.Range("F2").Copy($ws.Range("F5"))
$ws.Range("F5").Value = "Custom Potion"

# ------------------------------------------------------------------
# 7. Column widths for the new / shifted columns
# ------------------------------------------------------------------
$ws.Columns("B").ColumnWidth = 24.7109375
$ws.Columns("C").ColumnWidth = 52

# ------------------------------------------------------------------
# 8. Move the conditional formatting from column B to column D
# ------------------------------------------------------------------
$ws.Range("B1:B1048576").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D1:D1048576"))

# ------------------------------------------------------------------
# 9. Re-create the hyperlinks on their new cells
# ------------------------------------------------------------------
$bbmodelLink = "..\texturepack\assets\minecraft\models\item\custom\backpack.bbmodel"
$textureLink1 = "..\texturepack\assets\blockgame\textures\backpack.png"
$textureLink2 = "..\texturepack\assets\minecraft\models\item\bone_animated\bone_animated.png"
$jsonLink = "..\texturepack\assets\blockgame\models\bone_animated\bone_animated.json"

$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("G2"), $bbmodelLink, "", "", $bbmodelLink) | Out-Null
$ws.Hyperlinks.Add($ws.Range("H2"), $textureLink1, "", "", $textureLink1) | Out-Null
$ws.Hyperlinks.Add($ws.Range("H3"), $textureLink2, "", "", $textureLink2) | Out-Null
$ws.Hyperlinks.Add($ws.Range("G3"), $jsonLink, "", "", "") | Out-Null

# ------------------------------------------------------------------
# 10. Selection shown when the workbook is opened
# ------------------------------------------------------------------
$ws.Range("C8").Select()
